{"js": "// Change the evaluation scale numbering so that 5 is the best score and\n// 0 is the worst (was 1..6, now 5,4,3,2,1,0 respectively). Only the six\n// scale-list paragraphs (indented, \"N. <description>\") are touched; the\n// later \"Question\" list (which also starts with \"1.\", \"2.\", ...) must be\n// left untouched.\n\n// Map from the distinctive remainder of each scale-list paragraph (the\n// text right after \"N.\") to the old leading digit and the new digit that\n// should replace it.\nconst scaleItems = [\n  { match: \"All requirements were met and additional work was done to\", oldNum: \"1\", newNum: \"5\" },\n  { match: \"All requirements were met.\", oldNum: \"2\", newNum: \"4\" },\n  { match: \"Some attempt was made, but was significantly deficient\", oldNum: \"4\", newNum: \"2\" },\n  { match: \"Some attempt was made, but was extremely deficient\", oldNum: \"5\", newNum: \"1\" },\n  { match: \"No attempt was made.\", oldNum: \"6\", newNum: \"0\" }\n  // \"3. Some attempt was made, but was slightly deficient...\" is unchanged\n  // (3 stays 3), so it is intentionally not listed here.\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const item of paragraphs.items) {\n  const text = item.text;\n  for (const { match, oldNum, newNum } of scaleItems) {\n    const prefix = oldNum + \".\";\n    if (text.indexOf(match) !== -1 && text.indexOf(prefix) === 0) {\n      // Scope the search to this paragraph only, and match just the\n      // leading digit so the rest of the run (\".\" + description) is left\n      // completely untouched.\n      const found = item.search(oldNum, { matchCase: true, matchWholeWord: false });\n      found.load(\"items/text\");\n      await context.sync();\n\n      if (found.items.length > 0) {\n        found.items[0].insertText(newNum, \"Replace\");\n        await context.sync();\n      }\n      break;\n    }\n  }\n}\n", "ps1": "# Change the evaluation scale numbering so that 5 is the best score and\n# 0 is the worst (was 1..6, now 5,4,3,2,1,0 respectively). Only the six\n# scale-list paragraphs (indented, \"N. <description>\") are touched; the\n# later \"Question\" list (which also starts with \"1.\", \"2.\", ...) must be\n# left untouched.\n\n$d = $word.ActiveDocument\n\n# Map the distinctive remainder of each scale-list paragraph (the text\n# right after \"N.\") to the old leading digit and the new digit that\n# should replace it. \"3. Some attempt was made, but was slightly\n# deficient...\" is intentionally omitted because 3 stays 3.\n$scaleItems = @(\n    @{ Match = \"All requirements were met and additional work was done to\"; OldNum = \"1\"; NewNum = \"5\" },\n    @{ Match = \"All requirements were met.\"; OldNum = \"2\"; NewNum = \"4\" },\n    @{ Match = \"Some attempt was made, but was significantly deficient\"; OldNum = \"4\"; NewNum = \"2\" },\n    @{ Match = \"Some attempt was made, but was extremely deficient\"; OldNum = \"5\"; NewNum = \"1\" },\n    @{ Match = \"No attempt was made.\"; OldNum = \"6\"; NewNum = \"0\" }\n)\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    foreach ($item in $scaleItems) {\n        $prefix = $item.OldNum + \".\"\n        if ($text.StartsWith($prefix) -and $text.Contains($item.Match)) {\n            # Scope the Find/Replace to just this paragraph's range, and\n            # match only the single leading digit so the rest of the\n            # paragraph (the \".\" and description) is left untouched.\n            $rng = $p.Range\n            $find = $rng.Find\n            $find.ClearFormatting()\n            $find.Text = $item.OldNum\n            $find.Replacement.Text = $item.NewNum\n            $find.Forward = $true\n            $find.Wrap = 0\n            $find.MatchWildcards = $false\n            $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $item.NewNum, 2) | Out-Null\n            break\n        }\n    }\n}\n"}
